# Remove the second ("Cashflow") table that lived in columns M:P, rows 5-7.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cash Flows")

$ws.Range("M5:P7").Clear()

# Move the active selection to F15, matching the saved workbook state.
$ws.Range("F15").Select()
